$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in newly-measured "D" column (sheet column E) values
$ws.Range("E3").Value2 = -5.7
$ws.Range("E21").Value2 = -8.699999999999999

# These two measurements are now missing -> blank them out (stored as an
# empty text value, matching how the other missing cells in this sheet,
# e.g. C10/C12/C18, are represented)
$ws.Range("E5").Value2 = "'"
$ws.Range("E5").Style = "Normal"
$ws.Range("E23").Value2 = "'"
$ws.Range("E23").Style = "Normal"

# Remove the "RM 232" row entirely (row 26)
$ws.Rows(26).Delete()

# Remove what is now the "SC 92" row (shifted up into row 27)
$ws.Rows(27).Delete()

# Fill in the previously-missing "D" value for the "SC 193" row (now row 32)
$ws.Range("E32").Value2 = -6.4
